# Auto-generated edit script: updates cached profit-calculation values
# across the ALC/BSM/CRP/CUL/GSM/LTW/WVR sheets (per scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4217.811
$ws.Range("I11").Value = 4217.811
$ws.Range("K11").Value = 4217.811
$ws.Range("M11").Value = -4077.811
$ws.Range("H39").Value = 1055.1428
$ws.Range("I39").Value = 39.9
$ws.Range("K39").Value = 119.7
$ws.Range("M39").Value = 176.3
$ws.Range("H74").Value = 7057.125
$ws.Range("I74").Value = 6741.25
$ws.Range("K74").Value = 6741.25
$ws.Range("M74").Value = -5805.25
$ws.Range("H77").Value = 7057.125
$ws.Range("I77").Value = 6741.25
$ws.Range("K77").Value = 33706.25
$ws.Range("M77").Value = -29026.25
$ws.Range("H86").Value = 2656.6667
$ws.Range("I86").Value = 3483.3333
$ws.Range("J86").Value = 1830
$ws.Range("K86").Value = 3483.3333
$ws.Range("L86").Value = 1830
$ws.Range("M86").Value = -2360.3333
$ws.Range("N86").Value = -4076
$ws.Range("H89").Value = 2656.6667
$ws.Range("I89").Value = 3483.3333
$ws.Range("J89").Value = 1830
$ws.Range("K89").Value = 17416.6665
$ws.Range("L89").Value = 9150
$ws.Range("M89").Value = -11800.6665
$ws.Range("N89").Value = -20382
$ws.Range("H92").Value = 1808.037
$ws.Range("I92").Value = 1327.3334
$ws.Range("J92").Value = 3490.5
$ws.Range("K92").Value = 1327.3334
$ws.Range("L92").Value = 3490.5
$ws.Range("M92").Value = -79.33339999999998
$ws.Range("N92").Value = -5986.5
$ws.Range("H99").Value = 481
$ws.Range("I99").Value = 499.85715
$ws.Range("K99").Value = 1499.57145
$ws.Range("M99").Value = -1.571449999999913
$ws.Range("H111").Value = 4375.2856
$ws.Range("I111").Value = 4375.2856
$ws.Range("K111").Value = 13125.8568
$ws.Range("M111").Value = -10058.8568
$ws.Range("H125").Value = 4500
$ws.Range("I125").Value = 3714.2856
$ws.Range("J125").Value = 5600
$ws.Range("K125").Value = 33428.5704
$ws.Range("L125").Value = 50400
$ws.Range("M125").Value = -30968.5704
$ws.Range("N125").Value = -55320
$ws.Range("H138").Value = 3139.7222
$ws.Range("I138").Value = 2837.7856
$ws.Range("J138").Value = 4196.5
$ws.Range("K138").Value = 8513.356800000001
$ws.Range("L138").Value = 12589.5
$ws.Range("M138").Value = -3373.356800000001
$ws.Range("N138").Value = -22869.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1733.3
$ws.Range("I107").Value = 1324.6364
$ws.Range("K107").Value = 1324.6364
$ws.Range("M107").Value = 595.3635999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 23847.05
$ws.Range("J9").Value = 23847.05
$ws.Range("L9").Value = 23847.05
$ws.Range("N9").Value = -24183.05
$ws.Range("H86").Value = 9112.6
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 9112.6
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 9112.6
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -11358.6
$ws.Range("H89").Value = 9112.6
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 9112.6
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 45563
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -56795
$ws.Range("H105").Value = 2209.6667
$ws.Range("I105").Value = 1603.75
$ws.Range("K105").Value = 1603.75
$ws.Range("M105").Value = 143.25
$ws.Range("H107").Value = 45454868
$ws.Range("I107").Value = 52631812
$ws.Range("K107").Value = 52631812
$ws.Range("M107").Value = -52629892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4551762
$ws.Range("I4").Value = 3768559
$ws.Range("K4").Value = 11305677
$ws.Range("M4").Value = -11305565
$ws.Range("H137").Value = 6690.2915
$ws.Range("J137").Value = 3976.7693
$ws.Range("L137").Value = 11930.3079
$ws.Range("N137").Value = -22130.3079
$ws.Range("H140").Value = 2000.5385
$ws.Range("I140").Value = 2666
$ws.Range("J140").Value = 1945.0834
$ws.Range("K140").Value = 7998
$ws.Range("L140").Value = 5835.2502
$ws.Range("M140").Value = -2818
$ws.Range("N140").Value = -16195.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 88999
$ws.Range("J130").Value = 88999
$ws.Range("L130").Value = 88999
$ws.Range("N130").Value = -99039
$ws.Range("H132").Value = 5743.769
$ws.Range("I132").Value = 4018.6365
$ws.Range("J132").Value = 7976.294
$ws.Range("K132").Value = 12055.9095
$ws.Range("L132").Value = 23928.882
$ws.Range("M132").Value = -9525.9095
$ws.Range("N132").Value = -28988.882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 202400.8
$ws.Range("J7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("N7").Value = -2724
$ws.Range("H22").Value = 3315.175
$ws.Range("I22").Value = 2095.7368
$ws.Range("K22").Value = 2095.7368
$ws.Range("M22").Value = -1800.7368
$ws.Range("H27").Value = 3315.175
$ws.Range("I27").Value = 2095.7368
$ws.Range("K27").Value = 2095.7368
$ws.Range("M27").Value = -1988.7368
$ws.Range("H93").Value = 1170.069
$ws.Range("I93").Value = 1078.3684
$ws.Range("J93").Value = 1344.3
$ws.Range("K93").Value = 1078.3684
$ws.Range("L93").Value = 1344.3
$ws.Range("M93").Value = 169.6315999999999
$ws.Range("N93").Value = -3840.3
$ws.Range("H122").Value = 3073.125
$ws.Range("I122").Value = 3073.125
$ws.Range("K122").Value = 9219.375
$ws.Range("M122").Value = -6769.375
$ws.Range("H124").Value = 65265.4
$ws.Range("J124").Value = 65265.4
$ws.Range("L124").Value = 65265.4
$ws.Range("N124").Value = -75085.39999999999
$ws.Range("H126").Value = 202400.8
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440
$ws.Range("H134").Value = 59714.5
$ws.Range("J134").Value = 90429
$ws.Range("L134").Value = 90429
$ws.Range("N134").Value = -100569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 43960
$ws.Range("J120").Value = 43960
$ws.Range("L120").Value = 43960
$ws.Range("N120").Value = -53636
$ws.Range("H122").Value = 2231.1396
$ws.Range("I122").Value = 2081.6858
$ws.Range("K122").Value = 6245.057400000001
$ws.Range("M122").Value = -3795.057400000001
$ws.Range("H125").Value = 50715
$ws.Range("J125").Value = 50715
$ws.Range("L125").Value = 50715
$ws.Range("M125").Value = -60555
$ws.Range("H136").Value = 4929171.5
$ws.Range("I136").Value = 8405077
$ws.Range("J136").Value = 4971.25
$ws.Range("K136").Value = 25215231
$ws.Range("L136").Value = 14913.75
$ws.Range("M136").Value = -25212681
$ws.Range("N136").Value = -20013.75
